$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 130
$ws.Range("B130").Value = 7483189
$ws.Range("E130").Value = "Independiente del Valle"
$ws.Range("F130").Value = "Orense"
$ws.Range("G130").Value = 2
$ws.Range("H130").Value = 2
$ws.Range("J130").Value = 1.4
$ws.Range("K130").Value = 4.75
$ws.Range("L130").Value = 7
$ws.Range("M130").Value = 1.4
$ws.Range("N130").Value = 4.5
$ws.Range("O130").Value = 8
$ws.Range("P130").Value = -1.25
$ws.Range("Q130").Value = 1.875
$ws.Range("R130").Value = 1.925
$ws.Range("T130").Value = 1.925
$ws.Range("U130").Value = 1.875
$ws.Range("W130").Value = 3.5
$ws.Range("Z130").Value = 0.925
$ws.Range("AA130").Value = 0.925
$ws.Range("AB130").Value = -1
# Row 131
$ws.Range("B131").Value = 7483281
$ws.Range("E131").Value = "SD Aucas"
$ws.Range("F131").Value = "Delfin SC"
$ws.Range("G131").Value = 0
$ws.Range("I131").Value = "D"
$ws.Range("J131").Value = 1.909
$ws.Range("L131").Value = 4.2
$ws.Range("M131").Value = 1.909
$ws.Range("N131").Value = 3.5
$ws.Range("O131").Value = 4
$ws.Range("P131").Value = -0.5
$ws.Range("Q131").Value = 1.9
$ws.Range("R131").Value = 1.9
$ws.Range("S131").Value = 2.5
$ws.Range("T131").Value = 1.8
$ws.Range("U131").Value = 2
$ws.Range("V131").Value = -1
$ws.Range("W131").Value = 2.5
$ws.Range("Y131").Value = -1
$ws.Range("Z131").Value = 0.8999999999999999
$ws.Range("AB131").Value = 1
# Row 132
$ws.Range("B132").Value = 7483081
$ws.Range("E132").Value = "Deportivo Cuenca"
$ws.Range("F132").Value = "El Nacional"
$ws.Range("G132").Value = 1
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = "H"
$ws.Range("J132").Value = 2.75
$ws.Range("K132").Value = 3.25
$ws.Range("L132").Value = 2.55
$ws.Range("M132").Value = 3
$ws.Range("N132").Value = 3.3
$ws.Range("O132").Value = 2.3
$ws.Range("P132").Value = 0.25
$ws.Range("Q132").Value = 1.825
$ws.Range("R132").Value = 1.975
$ws.Range("S132").Value = 2.75
$ws.Range("T132").Value = 2
$ws.Range("U132").Value = 1.8
$ws.Range("V132").Value = 2
$ws.Range("X132").Value = -1
$ws.Range("Y132").Value = 0.825
$ws.Range("Z132").Value = -1
$ws.Range("AB132").Value = 0.8
# Row 133
$ws.Range("B133").Value = 7483247
$ws.Range("E133").Value = "Mushuc Runa"
$ws.Range("F133").Value = "Universidad Catolica del Ecuador"
$ws.Range("G133").Value = 0
$ws.Range("I133").Value = "A"
$ws.Range("J133").Value = 3.25
$ws.Range("K133").Value = 3.2
$ws.Range("L133").Value = 2.25
$ws.Range("M133").Value = 3.5
$ws.Range("N133").Value = 3.25
$ws.Range("O133").Value = 2.1
$ws.Range("P133").Value = 0.5
$ws.Range("Q133").Value = 1.775
$ws.Range("R133").Value = 2.025
$ws.Range("T133").Value = 1.9
$ws.Range("U133").Value = 1.9
$ws.Range("W133").Value = -1
$ws.Range("X133").Value = 1.1
$ws.Range("Z133").Value = 1.025
$ws.Range("AA133").Value = -1
$ws.Range("AB133").Value = 0.8999999999999999
# Row 134
$ws.Range("B134").Value = 7482832
$ws.Range("E134").Value = "Barcelona Guayaquil"
$ws.Range("F134").Value = "Guayaquil City"
$ws.Range("G134").Value = 2
$ws.Range("H134").Value = 1
$ws.Range("I134").Value = "H"
$ws.Range("J134").Value = 1.363
$ws.Range("K134").Value = 5
$ws.Range("L134").Value = 7.5
$ws.Range("M134").Value = 1.444
$ws.Range("N134").Value = 4
$ws.Range("O134").Value = 8
$ws.Range("P134").Value = -1.25
$ws.Range("Q134").Value = 2.05
$ws.Range("R134").Value = 1.75
$ws.Range("T134").Value = 1.95
$ws.Range("U134").Value = 1.85
$ws.Range("V134").Value = 0.444
$ws.Range("X134").Value = -1
$ws.Range("Y134").Value = -0.5
$ws.Range("Z134").Value = 0.375
$ws.Range("AA134").Value = 0.95
$ws.Range("AB134").Value = -1
# Row 137
$ws.Range("B137").Value = 7483188
$ws.Range("E137").Value = "Gualaceo SC"
$ws.Range("F137").Value = "Emelec"
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 2
$ws.Range("I137").Value = "A"
$ws.Range("J137").Value = 3.6
$ws.Range("K137").Value = 3.3
$ws.Range("L137").Value = 2.05
$ws.Range("M137").Value = 2.6
$ws.Range("N137").Value = 3.25
$ws.Range("O137").Value = 2.75
$ws.Range("P137").Value = 0
$ws.Range("Q137").Value = 1.8
$ws.Range("R137").Value = 2
$ws.Range("T137").Value = 1.975
$ws.Range("U137").Value = 1.825
$ws.Range("V137").Value = -1
$ws.Range("X137").Value = 1.75
$ws.Range("Y137").Value = -1
$ws.Range("Z137").Value = 1
$ws.Range("AA137").Value = -1
$ws.Range("AB137").Value = 0.825
# Row 142
$ws.Range("B142").Value = 7528857
$ws.Range("E142").Value = "Universidad Catolica del Ecuador"
$ws.Range("F142").Value = "Barcelona Guayaquil"
$ws.Range("G142").Value = 0
$ws.Range("I142").Value = "A"
$ws.Range("J142").Value = 1.533
$ws.Range("K142").Value = 4
$ws.Range("L142").Value = 5.5
$ws.Range("M142").Value = 1.5
$ws.Range("N142").Value = 4.333
$ws.Range("O142").Value = 5.25
$ws.Range("P142").Value = -1
$ws.Range("Q142").Value = 1.8
$ws.Range("R142").Value = 2
$ws.Range("S142").Value = 3
$ws.Range("T142").Value = 1.975
$ws.Range("U142").Value = 1.825
$ws.Range("V142").Value = -1
$ws.Range("X142").Value = 4.25
$ws.Range("Y142").Value = -1
$ws.Range("Z142").Value = 1
$ws.Range("AA142").Value = -1
$ws.Range("AB142").Value = 0.825
# Row 143
$ws.Range("B143").Value = 7528848
$ws.Range("E143").Value = "Emelec"
$ws.Range("F143").Value = "Deportivo Cuenca"
$ws.Range("G143").Value = 2
$ws.Range("H143").Value = 1
$ws.Range("I143").Value = "H"
$ws.Range("J143").Value = 1.75
$ws.Range("K143").Value = 3.5
$ws.Range("L143").Value = 4.2
$ws.Range("M143").Value = 2.4
$ws.Range("N143").Value = 3.1
$ws.Range("O143").Value = 2.75
$ws.Range("P143").Value = -0.25
$ws.Range("Q143").Value = 2.05
$ws.Range("R143").Value = 1.75
$ws.Range("T143").Value = 1.8
$ws.Range("U143").Value = 2
$ws.Range("V143").Value = 1.4
$ws.Range("X143").Value = -1
$ws.Range("Y143").Value = 1.05
$ws.Range("Z143").Value = -1
$ws.Range("AA143").Value = 0.8
# Row 145
$ws.Range("B145").Value = 7528858
$ws.Range("E145").Value = "Orense"
$ws.Range("F145").Value = "SD Aucas"
$ws.Range("G145").Value = 1
$ws.Range("H145").Value = 2
$ws.Range("J145").Value = 2.2
$ws.Range("K145").Value = 3.2
$ws.Range("L145").Value = 3.2
$ws.Range("M145").Value = 1.95
$ws.Range("N145").Value = 3.2
$ws.Range("O145").Value = 3.8
$ws.Range("P145").Value = -0.5
$ws.Range("Q145").Value = 1.95
$ws.Range("R145").Value = 1.85
$ws.Range("S145").Value = 2.25
$ws.Range("T145").Value = 1.85
$ws.Range("U145").Value = 1.95
$ws.Range("X145").Value = 2.8
$ws.Range("Z145").Value = 0.8500000000000001
$ws.Range("AA145").Value = 0.8500000000000001
$ws.Range("AB145").Value = -1
# Row 216
$ws.Range("M216").Value = 1.285
$ws.Range("O216").Value = 7.5
$ws.Range("P216").Value = -1.5
$ws.Range("Q216").Value = 1.8
$ws.Range("R216").Value = 2
$ws.Range("S216").Value = 3
$ws.Range("T216").Value = 1.925
$ws.Range("U216").Value = 1.875
# Row 217
$ws.Range("M217").Value = 2.45
$ws.Range("N217").Value = 2.875
$ws.Range("O217").Value = 3.2
$ws.Range("S217").Value = 2
$ws.Range("T217").Value = 1.95
$ws.Range("U217").Value = 1.85
# Row 218
$ws.Range("M218").Value = 1.5
$ws.Range("N218").Value = 4
$ws.Range("O218").Value = 6
$ws.Range("T218").Value = 1.85
$ws.Range("U218").Value = 1.95
# Row 219
$ws.Range("M219").Value = 4.75
$ws.Range("O219").Value = 1.65
$ws.Range("Q219").Value = 1.95
$ws.Range("R219").Value = 1.85
# Row 220
$ws.Range("M220").Value = 2.6
$ws.Range("N220").Value = 3.2
$ws.Range("O220").Value = 2.6
$ws.Range("Q220").Value = 1.9
$ws.Range("R220").Value = 1.9
$ws.Range("S220").Value = 2.25
$ws.Range("T220").Value = 1.775
$ws.Range("U220").Value = 2.025
# Row 221
$ws.Range("M221").Value = 2.25
$ws.Range("O221").Value = 3.2
$ws.Range("T221").Value = 1.85
$ws.Range("U221").Value = 1.95
